$d = $word.ActiveDocument

# --- 1. Create the new MSC_Join_A / MSC_Join_B / MSC_Join_C paragraph
#        styles, each based on the existing MSCJoin style. ---

$joinA = $d.Styles.Add("MSC_Join_A", 1)
$joinA.BaseStyle = $d.Styles("MSCJoin")

$joinB = $d.Styles.Add("MSC_Join_B", 1)
$joinB.BaseStyle = $d.Styles("MSCJoin")
$joinB.Font.Name = "Noto Sans CJK SC"
$joinB.Font.NameFarEast = "Noto Sans CJK SC"
$joinB.Font.NameBi = "Noto Sans CJK SC"

$joinC = $d.Styles.Add("MSC_Join_C", 1)
$joinC.BaseStyle = $d.Styles("MSCJoin")

# --- 2. Re-point every paragraph that currently uses the generic
#        "MSC_Join" style to the variant that matches the column
#        (A/B/C) it sits in, tracked via the most recent
#        MSC_Paragraph_A / MSC_Paragraph_B / MSC_Paragraph_C
#        paragraph seen while walking the document. ---

$currentVariant = $null

foreach ($p in $d.Paragraphs) {
    $styleName = $p.Style.NameLocal

    if ($styleName -eq "MSC_Paragraph_A") {
        $currentVariant = "A"
    }
    elseif ($styleName -eq "MSC_Paragraph_B") {
        $currentVariant = "B"
    }
    elseif ($styleName -eq "MSC_Paragraph_C") {
        $currentVariant = "C"
    }
    elseif ($styleName -eq "MSC_Join") {
        if ($currentVariant -eq "A") {
            $p.Style = "MSC_Join_A"
        }
        elseif ($currentVariant -eq "B") {
            $p.Style = "MSC_Join_B"
        }
        elseif ($currentVariant -eq "C") {
            $p.Style = "MSC_Join_C"
        }
    }
}
